function Find-ParagraphIndex($doc, $substr) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($substr)) {
            return $i
        }
    }
    return -1
}

$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1. Title line: merge "Springboard Data Science" + bookmark + " Intensive"
#    into a single run "Springboard Data Science Intensive" (the _GoBack
#    bookmark moves further down the document, see step 4).
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "Springboard Data Science"
$p = $d.Paragraphs.Item($idx)
$xml = "<w:p $wns><w:pPr><w:spacing w:line='480' w:lineRule='auto'/><w:jc w:val='right'/></w:pPr>" + `
       "<w:r><w:t>Springboard Data Science Intensive</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2. Replace the paragraph that begins "I thought there may be some merit..."
#    with the new "first thing I did" (random forest) paragraph text, and
#    split the remaining "What really interested me..." sentence into its
#    own new paragraph right after it.
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "I thought there may be some merit"
$p = $d.Paragraphs.Item($idx)
$xml = "<w:p $wns><w:pPr><w:spacing w:line='480' w:lineRule='auto'/></w:pPr>" + `
       "<w:r><w:tab/><w:t xml:space='preserve'>The first thing I did was a quick analysis to determine if who is consuming energy has changed enough over the decades that a forest could find it. I broke up a subset of dataframe that contained information about what industries are consuming the generated energy. A quick graph showed us that the mix had been changing over the years. Using a random forest, I ran a grid search cross validation and found my </w:t></w:r>" + `
       "<w:proofErr w:type='spellStart'/><w:r><w:t>n_estimators</w:t></w:r><w:proofErr w:type='spellEnd'/>" + `
       "<w:r><w:t xml:space='preserve'> should be three and I should not set a </w:t></w:r>" + `
       "<w:proofErr w:type='spellStart'/><w:r><w:t>max_depth</w:t></w:r><w:proofErr w:type='spellEnd'/>" + `
       "<w:r><w:t>. This forest predicted our decade based on the breakdown of who was using power in a given year, and yielded a ~88% accuracy.</w:t></w:r>" + `
       "</w:p>"
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$newP = $d.Paragraphs.Item($idx + 1)
$xml = "<w:p $wns><w:pPr><w:spacing w:line='480' w:lineRule='auto'/></w:pPr>" + `
       "<w:r><w:tab/></w:r>" + `
       "<w:r><w:t>What really interested me was whether or not states are using coal. The reason I’m picking on coal is because it’s the only dirty fuel still pervasive in today’s energy mix. Natural gas is a fossil fuel, but its emissions are so small (relative to coal), that even most progressives are accepting of it for the near-mid future.</w:t></w:r>" + `
       "</w:p>"
$newP.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3. "The first step in this data frame was to group by state..." gains an
#    inserted "'s main use" clause.
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "The first step in this data frame was to group by state"
$p = $d.Paragraphs.Item($idx)
$xml = "<w:p $wns><w:pPr><w:spacing w:line='480' w:lineRule='auto'/></w:pPr>" + `
       "<w:r><w:tab/><w:t>The first step in this data frame</w:t></w:r>" + `
       "<w:r><w:t>’s main use</w:t></w:r>" + `
       "<w:r><w:t xml:space='preserve'> was to group by state and input type. From here, I sorted and found the max generation associated with a given input type for each state. This was a big assumption on my part, both because this was a sum of all years, meaning states that had recently switched to cleaner energy may still be categorized in the coal category, and it didn’t include the second biggest generation. Regardless, I moved forward because I needed some form of analysis, and those would certainly do. The dataframe I performed analysis on was a surprisingly good representation of what is going on, even with all the assumptions.</w:t></w:r>" + `
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 4. Insert a new paragraph right after "Wanting to be thorough I did an
#    A/B test..." describing a second A/B test, carrying the _GoBack
#    bookmark with it.
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "Wanting to be thorough I did an A/B test"
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$newP = $d.Paragraphs.Item($idx + 1)
$xml = "<w:p $wns><w:pPr><w:spacing w:line='480' w:lineRule='auto'/></w:pPr>" + `
       "<w:r><w:tab/><w:t>I did another A/B test testing whether a Random Forest could do a better job than our K-nearest neighbors on predicting based on quantity of generation alone. I found that there is a significant</w:t></w:r>" + `
       "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" + `
       "<w:r><w:t xml:space='preserve'> difference between the two, with the Random Forest doing better, at the 85% level.</w:t></w:r>" + `
       "</w:p>"
$newP.Range.InsertXML($xml)
